# Regenerate save data: column G ("K") is recomputed (replacing the old
# "Strike#" derived value) and rewritten for each data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values keyed by row number (row 1 is the header row; data starts at row 2)
$kValues = @{
    2  = 0
    3  = 1
    4  = 2
    5  = 1
    6  = 1
    7  = 2
    8  = 3
    9  = 0
    10 = 0
    11 = 3
    12 = 3
    13 = 1
    14 = 0
    15 = 0
    16 = 2
    17 = 0
    18 = 0
    19 = 1
    20 = 0
    21 = 0
    22 = 2
    23 = 1
    24 = 2
    25 = 0
    26 = 3
    27 = 3
    28 = 0
    29 = 3
    30 = 2
    31 = 2
    32 = 1
    33 = 1
    34 = 0
    35 = 1
    36 = 0
    37 = 1
    38 = 1
    39 = 2
    40 = 0
    41 = 1
    42 = 1
    43 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
